$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns to English snake_case names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case municipality / state name fixes (and the MonteMorelos -> Montemorelos fix) ---
$ws.Range("B8").Value = "Pabellón De Arteaga"
$ws.Range("B9").Value = "Rincón De Romos"
$ws.Range("B36").Value = "Amatenango De La Frontera"
$ws.Range("B37").Value = "Amatenango Del Valle"
$ws.Range("B40").Value = "Bejucal De Ocampo"
$ws.Range("B50").Value = "Chiapa De Corzo"
$ws.Range("B57").Value = "Comitán De Domínguez"
$ws.Range("B86").Value = "Marqués De Comillas"
$ws.Range("B87").Value = "Mazapa De Madero"
$ws.Range("B91").Value = "Montecristo De Guerrero"
$ws.Range("B96").Value = "Ocozocoautla De Espinosa"
$ws.Range("B108").Value = "Salto De Agua"
$ws.Range("B109").Value = "San Cristóbal De Las Casas"
$ws.Range("B113").Value = "Santiago El Pinar"
$ws.Range("B150").Value = "Guadalupe Y Calvo"
$ws.Range("B152").Value = "Hidalgo Del Parral"
$ws.Range("B159").Value = "San Francisco Del Oro"
$ws.Range("B161").Value = "Valle De Zaragoza"
$ws.Range("B178").Value = "San Juan De Sabinas"
$ws.Range("B189").Value = "Villa De Álvarez"
$ws.Range("A191").Value = "Ciudad De México"
$ws.Range("B195").Value = "Cuajimalpa De Morelos"
$ws.Range("B222").Value = "Nombre De Dios"
$ws.Range("B231").Value = "San Juan De Guadalupe"
$ws.Range("B232").Value = "San Juan Del Río"
$ws.Range("A240").Value = "Estado De México"
$ws.Range("B240").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B242").Value = "Almoloya De Alquisiras"
$ws.Range("B243").Value = "Almoloya De Juárez"
$ws.Range("B244").Value = "Almoloya Del Río"
$ws.Range("B250").Value = "Atizapán De Zaragoza"
$ws.Range("B258").Value = "Coacalco De Berriozábal"
$ws.Range("B263").Value = "Ecatepec De Morelos"
$ws.Range("B270").Value = "Ixtapan De La Sal"
$ws.Range("B283").Value = "Naucalpan De Juárez"
$ws.Range("B294").Value = "San Antonio La Isla"
$ws.Range("B295").Value = "San Felipe Del Progreso"
$ws.Range("B297").Value = "San Simón De Guerrero"
$ws.Range("B307").Value = "Tenango Del Valle"
$ws.Range("B317").Value = "Tlalnepantla De Baz"
$ws.Range("B322").Value = "Valle De Bravo"
$ws.Range("B323").Value = "Valle De Chalco Solidaridad"
$ws.Range("B324").Value = "Villa De Allende"
$ws.Range("B325").Value = "Villa Del Carbón"
$ws.Range("B338").Value = "Apaseo El Alto"
$ws.Range("B339").Value = "Apaseo El Grande"
$ws.Range("B347").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B351").Value = "Jaral Del Progreso"
$ws.Range("B359").Value = "Purísima Del Rincón"
$ws.Range("B363").Value = "San Diego De La Unión"
$ws.Range("B365").Value = "San Francisco Del Rincón"
$ws.Range("B367").Value = "San Luis De La Paz"
$ws.Range("B369").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B371").Value = "Silao De La Victoria"
$ws.Range("B376").Value = "Valle De Santiago"
$ws.Range("B382").Value = "Acapulco De Juárez"
$ws.Range("B385").Value = "Ajuchitlán Del Progreso"
$ws.Range("B386").Value = "Alcozauca De Guerrero"
$ws.Range("B390").Value = "Atenango Del Río"
$ws.Range("B391").Value = "Atlamajalcingo Del Monte"
$ws.Range("B393").Value = "Atoyac De Álvarez"
$ws.Range("B394").Value = "Ayutla De Los Libres"
$ws.Range("B397").Value = "Buenavista De Cuéllar"
$ws.Range("B398").Value = "Chilapa De Álvarez"
$ws.Range("B399").Value = "Chilpancingo De Los Bravo"
$ws.Range("B400").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B405").Value = "Coyuca De Benítez"
$ws.Range("B406").Value = "Coyuca De Catalán"
$ws.Range("B410").Value = "Cuetzala Del Progreso"
$ws.Range("B411").Value = "Cutzamala De Pinzón"
$ws.Range("B417").Value = "Huitzuco De Los Figueroa"
$ws.Range("B418").Value = "Iguala De La Independencia"
$ws.Range("B420").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B421").Value = "Zihuatanejo De Azueta"
$ws.Range("B423").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B426").Value = "Mártir De Cuilapan"
$ws.Range("B439").Value = "Taxco De Alarcón"
$ws.Range("B441").Value = "Técpan De Galeana"
$ws.Range("B443").Value = "Tepecoacuilco De Trujano"
$ws.Range("B445").Value = "Tixtla De Guerrero"
$ws.Range("B449").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B450").Value = "Tlapa De Comonfort"
$ws.Range("B462").Value = "Agua Blanca De Iturbide"
$ws.Range("B468").Value = "Atotonilco De Tula"
$ws.Range("B469").Value = "Atotonilco El Grande"
$ws.Range("B475").Value = "Cuautepec De Hinojosa"
$ws.Range("B480").Value = "Huasca De Ocampo"
$ws.Range("B484").Value = "Huejutla De Reyes"
$ws.Range("B487").Value = "Jacala De Ledezma"
$ws.Range("B493").Value = "Mineral Del Chico"
$ws.Range("B494").Value = "Mineral Del Monte"
$ws.Range("B495").Value = "Mixquiahuala De Juárez"
$ws.Range("B496").Value = "Molango De Escamilla"
$ws.Range("B498").Value = "Nopala De Villagrán"
$ws.Range("B499").Value = "Pachuca De Soto"
$ws.Range("B502").Value = "Progreso De Obregón"
$ws.Range("B508").Value = "Santiago De Anaya"
$ws.Range("B509").Value = "Santiago Tulantepec De Lugo Guerrero"
$ws.Range("B513").Value = "Tenango De Doria"
$ws.Range("B515").Value = "Tepehuacán De Guerrero"
$ws.Range("B516").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B517").Value = "Tezontepec De Aldama"
$ws.Range("B525").Value = "Tula De Allende"
$ws.Range("B526").Value = "Tulancingo De Bravo"
$ws.Range("B527").Value = "Villa De Tezontepec"
$ws.Range("B531").Value = "Zacualtipán De Ángeles"
$ws.Range("B536").Value = "Ahualulco De Mercado"
$ws.Range("B541").Value = "Atotonilco El Alto"
$ws.Range("B543").Value = "Autlán De Navarro"
$ws.Range("B548").Value = "Cañadas De Obregón"
$ws.Range("B556").Value = "Encarnación De Díaz"
$ws.Range("B560").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B561").Value = "Ixtlahuacán Del Río"
$ws.Range("B564").Value = "Jilotlán De Los Dolores"
$ws.Range("B569").Value = "Lagos De Moreno"
$ws.Range("B576").Value = "Ojuelos De Jalisco"
$ws.Range("B581").Value = "San Cristóbal De La Barranca"
$ws.Range("B582").Value = "San Diego De Alejandría"
$ws.Range("B583").Value = "San Juan De Los Lagos"
$ws.Range("B585").Value = "San Miguel El Alto"
$ws.Range("B586").Value = "San Sebastián Del Oeste"
$ws.Range("B589").Value = "Tamazula De Gordiano"
$ws.Range("B592").Value = "Teocuitatlán De Corona"
$ws.Range("B593").Value = "Tepatitlán De Morelos"
$ws.Range("B595").Value = "Tizapán El Alto"
$ws.Range("B596").Value = "Tlajomulco De Zúñiga"
$ws.Range("B602").Value = "Unión De San Antonio"
$ws.Range("B603").Value = "Unión De Tula"
$ws.Range("B604").Value = "Valle De Juárez"
$ws.Range("B607").Value = "Yahualica De González Gallo"
$ws.Range("B608").Value = "Zacoalco De Torres"
$ws.Range("B611").Value = "Zapotlán El Grande"
$ws.Range("B633").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B698").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B721").Value = "Coatlán Del Río"
$ws.Range("B733").Value = "Puente De Ixtla"
$ws.Range("B739").Value = "Tetela Del Volcán"
$ws.Range("B740").Value = "Tlaltizapán De Zapata"
$ws.Range("B746").Value = "Zacualpan De Amilpas"
$ws.Range("B750").Value = "Bahía De Banderas"
$ws.Range("B753").Value = "Ixtlán Del Río"
$ws.Range("B758").Value = "Santa María Del Oro"
$ws.Range("B776").Value = "Lampazos De Naranjo"
$ws.Range("B778").Value = "Mier Y Noriega"
$ws.Range("B783").Value = "San Nicolás De Los Garza"
$ws.Range("B787").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B796").Value = "Chalcatongo De Hidalgo"
$ws.Range("B797").Value = "Chiquihuitlán De Benito Juárez"
$ws.Range("B799").Value = "Coicoyán De Las Flores"
$ws.Range("B802").Value = "Cuilápam De Guerrero"
$ws.Range("B803").Value = "Cuyamecalco Villa De Zaragoza"
$ws.Range("B805").Value = "Fresnillo De Trujano"
$ws.Range("B806").Value = "Guadalupe De Ramírez"
$ws.Range("B807").Value = "Guevea De Humboldt"
$ws.Range("B808").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B809").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B810").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B811").Value = "Ixtlán De Juárez"
$ws.Range("B812").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B818").Value = "Mariscala De Juárez"
$ws.Range("B821").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B823").Value = "Nejapa De Madero"
$ws.Range("B825").Value = "Oaxaca De Juárez"
$ws.Range("B826").Value = "Ocotlán De Morelos"
$ws.Range("B827").Value = "Pinotepa De Don Luis"
$ws.Range("B829").Value = "Putla Villa De Guerrero"
$ws.Range("B830").Value = "Reforma De Pineda"
$ws.Range("B838").Value = "San Antonino El Alto"
$ws.Range("B850").Value = "San Dionisio Del Mar"
$ws.Range("B853").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B856").Value = "San Francisco Del Mar"
$ws.Range("B873").Value = "San José Del Progreso"
$ws.Range("B888").Value = "San Juan De Los Cués"
$ws.Range("B889").Value = "San Juan Del Estado"
$ws.Range("B890").Value = "San Juan Del Río"
$ws.Range("B931").Value = "San Miguel Del Puerto"
$ws.Range("B932").Value = "San Miguel El Grande"
$ws.Range("B943").Value = "San Pedro El Alto"
$ws.Range("B954").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B955").Value = "San Pedro Y San Pablo Tequixtepec"
$ws.Range("B978").Value = "Santa Cruz Tacache De Mina"
$ws.Range("B983").Value = "Santa Lucía Del Camino"
$ws.Range("B995").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B1047").Value = "Santo Domingo De Morelos"
$ws.Range("B1066").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B1067").Value = "Tanetze De Zaragoza"
$ws.Range("B1068").Value = "Tataltepec De Valdés"
$ws.Range("B1069").Value = "Teotitlán Del Valle"
$ws.Range("B1070").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B1071").Value = "Tlacolula De Matamoros"
$ws.Range("B1072").Value = "Totontepec Villa De Morelos"
$ws.Range("B1076").Value = "Villa De Chilapa De Díaz"
$ws.Range("B1077").Value = "Villa De Etla"
$ws.Range("B1078").Value = "Villa De Tamazulápam Del Progreso"
$ws.Range("B1079").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B1080").Value = "Villa Sola De Vega"
$ws.Range("B1081").Value = "Villa Tejúpam De La Unión"
$ws.Range("B1082").Value = "Zapotitlán Del Río"
$ws.Range("B1084").Value = "Zimatlán De Álvarez"
$ws.Range("B1101").Value = "Ayotoxco De Guerrero"
$ws.Range("B1112").Value = "Chila De La Sal"
$ws.Range("B1118").Value = "Cuapiaxtla De Madero"
$ws.Range("B1120").Value = "Cuayuca De Andrade"
$ws.Range("B1121").Value = "Cuetzalan Del Progreso"
$ws.Range("B1131").Value = "Huehuetlán El Chico"
$ws.Range("B1136").Value = "Izúcar De Matamoros"
$ws.Range("B1142").Value = "Los Reyes De Juárez"
$ws.Range("B1151").Value = "Palmar De Bravo"
$ws.Range("B1159").Value = "San Diego La Mesa Tochimiltzingo"
$ws.Range("B1165").Value = "San Nicolás De Los Ranchos"
$ws.Range("B1167").Value = "San Salvador El Seco"
$ws.Range("B1168").Value = "San Salvador El Verde"
$ws.Range("B1172").Value = "Tecali De Herrera"
$ws.Range("B1178").Value = "Tepanco De López"
$ws.Range("B1182").Value = "Tepexi De Rodríguez"
$ws.Range("B1184").Value = "Tetela De Ocampo"
$ws.Range("B1185").Value = "Teteles De Avila Castillo"
$ws.Range("B1189").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B1203").Value = "Xayacatlán De Bravo"
$ws.Range("B1207").Value = "Xochitlán De Vicente Suárez"
$ws.Range("B1220").Value = "Amealco De Bonfil"
$ws.Range("B1222").Value = "Cadereyta De Montes"
$ws.Range("B1228").Value = "Jalpan De Serra"
$ws.Range("B1229").Value = "Landa De Matamoros"
$ws.Range("B1232").Value = "Pinal De Amoles"
$ws.Range("B1235").Value = "San Juan Del Río"
$ws.Range("B1248").Value = "Armadillo De Los Infante"
$ws.Range("B1249").Value = "Axtla De Terrazas"
$ws.Range("B1255").Value = "Ciudad Del Maíz"
$ws.Range("B1265").Value = "Mexquitic De Carmona"
$ws.Range("B1271").Value = "San Ciro De Acosta"
$ws.Range("B1277").Value = "Santa María Del Río"
$ws.Range("B1279").Value = "Soledad De Graciano Sánchez"
$ws.Range("B1286").Value = "Tanquián De Escobedo"
$ws.Range("B1289").Value = "Villa De Arista"
$ws.Range("B1290").Value = "Villa De Arriaga"
$ws.Range("B1291").Value = "Villa De Guadalupe"
$ws.Range("B1292").Value = "Villa De La Paz"
$ws.Range("B1293").Value = "Villa De Ramos"
$ws.Range("B1294").Value = "Villa De Reyes"
$ws.Range("B1338").Value = "Jalpa De Méndez"
$ws.Range("B1377").Value = "Soto La Marina"
$ws.Range("B1387").Value = "Contla De Juan Cuamatzi"
$ws.Range("B1391").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B1393").Value = "Muñoz De Domingo Arenas"
$ws.Range("B1394").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B1399").Value = "San Pablo Del Monte"
$ws.Range("B1413").Value = "Ziltlaltépec De Trinidad Sánchez Santos"
$ws.Range("B1423").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B1426").Value = "Amatlán De Los Reyes"
$ws.Range("B1436").Value = "Boca Del Río"
$ws.Range("B1441").Value = "Castillo De Teayo"
$ws.Range("B1457").Value = "Cosamaloapan De Carpio"
$ws.Range("B1458").Value = "Cosautlán De Carvajal"
$ws.Range("B1475").Value = "Hueyapan De Ocampo"
$ws.Range("B1476").Value = "Ignacio De La Llave"
$ws.Range("B1479").Value = "Ixhuatlán De Madero"
$ws.Range("B1480").Value = "Ixhuatlán Del Café"
$ws.Range("B1481").Value = "Ixhuatlán Del Sureste"
$ws.Range("B1489").Value = "Juchique De Ferrer"
$ws.Range("B1493").Value = "Lerdo De Tejada"
$ws.Range("B1497").Value = "Martínez De La Torre"
$ws.Range("B1499").Value = "Medellín De Bravo"
$ws.Range("B1503").Value = "Mixtla De Altamirano"
$ws.Range("B1505").Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Range("B1517").Value = "Paso De Ovejas"
$ws.Range("B1518").Value = "Paso Del Macho"
$ws.Range("B1522").Value = "Poza Rica De Hidalgo"
$ws.Range("B1530").Value = "Sayula De Alemán"
$ws.Range("B1533").Value = "Soledad De Doblado"
$ws.Range("B1540").Value = "Tatahuicapan De Juárez"
$ws.Range("B1572").Value = "Vega De Alatorre"
$ws.Range("B1583").Value = "Zontecomatlán De López Y Fuentes"
$ws.Range("B1599").Value = "Cañitas De Felipe Pescador"
$ws.Range("B1601").Value = "Concepción Del Oro"
$ws.Range("B1603").Value = "El Plateado De Joaquín Amaro"
$ws.Range("B1611").Value = "Jiménez Del Teul"
$ws.Range("B1619").Value = "Nochistlán De Mejía"
$ws.Range("B1620").Value = "Noria De Ángeles"
$ws.Range("B1629").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B1633").Value = "Villa De Cos"
$ws.Range("B780").Value = "Montemorelos"

# --- Floating point literal (1-ULP) refresh on percentage column ---
$ws.Range("D43").Value = 0.0009539088583445344
$ws.Range("D69").Value = 0.0009539088583445344
$ws.Range("D155").Value = 0.0009539088583445344
$ws.Range("D174").Value = 0.0009539088583445344
$ws.Range("D259").Value = 0.0009539088583445344
$ws.Range("D288").Value = 0.0009539088583445344
$ws.Range("D339").Value = 0.0009539088583445344
$ws.Range("D418").Value = 0.0009972683519056495
$ws.Range("D473").Value = 0.0009972683519056495
$ws.Range("D692").Value = 0.0009972683519056495
$ws.Range("D809").Value = 0.0009539088583445344
$ws.Range("D819").Value = 0.0009972683519056495
$ws.Range("D925").Value = 0.0009539088583445344
$ws.Range("D992").Value = 0.0009539088583445344
$ws.Range("D1027").Value = 0.0009539088583445344
$ws.Range("D1085").Value = 0.09833933139660928
$ws.Range("D1109").Value = 0.0009972683519056495
$ws.Range("D1242").Value = 0.0009539088583445344
$ws.Range("D1305").Value = 0.0009539088583445344
$ws.Range("D1457").Value = 0.0009972683519056495
$ws.Range("D1472").Value = 0.0009972683519056495
$ws.Range("D1522").Value = 0.0009539088583445344
$ws.Range("D1543").Value = 0.0009539088583445344

# --- Remove trailing metadata/footer rows (1642-1646); row 1641 was already blank ---
$ws.Rows("1642:1646").Delete()
